$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'91.169.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +3.68%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.106.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.67%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.10%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'219.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +4.45%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'622.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.01%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.378"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +2.60%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.972"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +21.76%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.01%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'3.104.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.68%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.720"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +21.88%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +5.44%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +8.22%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'34.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +8.25%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "'Toncoin"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'5.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.37%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'WrappedBTC"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'90.997.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.59%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.689.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.69%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.138.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.51%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +13.97%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +10.57%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'14.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +5.96%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'434.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +3.66%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +8.41%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +6.45%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'6.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +13.22%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'87.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +6.62%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'12.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +3.54%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'3.284.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.64%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.10%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.168"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.60%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +13.49%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'524.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +4.15%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.891"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -17.62%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'3.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +4.97%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +5.55%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +10.33%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'23.66"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +6.78%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +4.18%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +3.96%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'WhiteBITCoin"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'22.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.37%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'Hedera"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.0861"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +25.25%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +0.01%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.152"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +15.28%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.391"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +8.88%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +0.00%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +6.51%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'146.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.59%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'44.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.53%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +10.23%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'166.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +6.24%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +7.05%  "
$ws.Range("E51").Style = "Normal"
